# Update the "Ressourcentyp, allgemein" / "Ressourcentyp, speziell" test
# placeholder values on the "Probe1" sheet to the real model values used
# for the biospecimen (Bioproben) resource type update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probe1")

# Make sure we're working on the correct sheet (it was already the active one)
$ws.Activate()

$ws.Range("B12").Value = "Bioproben"
$ws.Range("B13").Value = "Infektionskrankheiten"

# Leave the selection on the last edited cell, matching the authored file
$ws.Range("B13").Select()

$wb.Save()
